# Applies the "Add files via upload" edit to the "gUSD 26.06.25" sheet:
#  - fills in rows 143..149 for columns C,D,E,F,G (raw data), I (price formula),
#    M (sum formula), extending the shared formulas that previously stopped at row 142
#  - leaves rows 150..160 as-is for content, but Excel will naturally adjust the
#    <row> spans once nearby shared-formula ranges grow (handled by the engine)
#  - moves the active selection to G150

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("gUSD 26.06.25")
$ws.Activate()

# Raw input data for rows 143-149 (columns C,D,E,F,G)
$data = @(
    @{ Row = 143; C = 34.2109;  D = 95.9289;  E = 9.41; F = 12.26; G = 10.68 },
    @{ Row = 144; C = 32.8576;  D = 98.0636;  E = 9.42; F = 10.98; G = 16.49 },
    @{ Row = 145; C = 32.6985;  D = 98.5604;  E = 9.70; F = 11.36; G = 8.88  },
    @{ Row = 146; C = 31.1592;  D = 99.0279;  E = 9.68; F = 11.27; G = 4.50  },
    @{ Row = 147; C = 29.6177;  D = 99.5872;  E = 9.57; F = 9.75;  G = 4.20  },
    @{ Row = 148; C = 28.0210;  D = 100.667;  E = 9.46; F = 8.36;  G = 8.08  },
    @{ Row = 149; C = 26.7268;  D = 101.348;  E = 9.46; F = 7.37;  G = 5.07  }
)

foreach ($row in $data) {
    $r = $row.Row
    $ws.Cells.Item($r, 3).Value = $row.C   # C
    $ws.Cells.Item($r, 4).Value = $row.D   # D
    $ws.Cells.Item($r, 5).Value = $row.E   # E
    $ws.Cells.Item($r, 6).Value = $row.F   # F
    $ws.Cells.Item($r, 7).Value = $row.G   # G
}

# Extend the shared formulas for I (price) and M (sum) from row 142 down to 149,
# matching the existing formula pattern used for rows 65..142.
$ws.Range("I65:I149").Formula = "=C65/`$D`$3"
$ws.Range("M65:M149").Formula = "=C65+D65"

# The newly-written I143:I149 formula cells need the same number format as the
# rest of the I column (percentage-style "0.0000"), since writing .Formula on a
# previously-empty cell resets it to the default style.
$ws.Range("I143:I149").NumberFormat = "0.0000"

# Move the active selection/top-left cell to match the saved view state.
$ws.Range("G150").Select()
